$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 101
$ws.Range("I12").Value = 101
$ws.Range("K12").Value = 101
$ws.Range("M12").Value = 69
$ws.Range("H64").Value = 38546.25
$ws.Range("I64").Value = 127612.25
$ws.Range("J64").Value = 2919.85
$ws.Range("K64").Value = 127612.25
$ws.Range("L64").Value = 2919.85
$ws.Range("M64").Value = -127364.25
$ws.Range("N64").Value = -3415.85
$ws.Range("H67").Value = 38546.25
$ws.Range("I67").Value = 127612.25
$ws.Range("J67").Value = 2919.85
$ws.Range("K67").Value = 127612.25
$ws.Range("L67").Value = 2919.85
$ws.Range("M67").Value = -126754.25
$ws.Range("N67").Value = -4635.85
$ws.Range("H76").Value = 3066.8096
$ws.Range("I76").Value = 3000.2307
$ws.Range("J76").Value = 3175
$ws.Range("K76").Value = 3000.2307
$ws.Range("L76").Value = 3175
$ws.Range("M76").Value = -2685.2307
$ws.Range("N76").Value = -3805
$ws.Range("H79").Value = 3066.8096
$ws.Range("I79").Value = 3000.2307
$ws.Range("J79").Value = 3175
$ws.Range("K79").Value = 3000.2307
$ws.Range("L79").Value = 3175
$ws.Range("M79").Value = -1908.2307
$ws.Range("N79").Value = -5359
$ws.Range("H137").Value = 1236.8485
$ws.Range("I137").Value = 765.1142599999999
$ws.Range("J137").Value = 1494.8281
$ws.Range("K137").Value = 2295.34278
$ws.Range("L137").Value = 4484.4843
$ws.Range("M137").Value = 254.6572200000001
$ws.Range("N137").Value = -9584.4843

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28674.627
$ws.Range("I32").Value = 29888.541
$ws.Range("K32").Value = 29888.541
$ws.Range("M32").Value = -29601.541
$ws.Range("H53").Value = 16650
$ws.Range("J53").Value = 19975
$ws.Range("L53").Value = 19975
$ws.Range("N53").Value = -21339
$ws.Range("H113").Value = 40650.57
$ws.Range("J113").Value = 40650.57
$ws.Range("L113").Value = 40650.57
$ws.Range("N113").Value = -49328.57

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2272.5454
$ws.Range("I105").Value = 2185.625
$ws.Range("J105").Value = 2322.2144
$ws.Range("K105").Value = 2185.625
$ws.Range("L105").Value = 2322.2144
$ws.Range("M105").Value = -438.625
$ws.Range("N105").Value = -5816.2144
$ws.Range("H126").Value = 50776
$ws.Range("J126").Value = 50776
$ws.Range("L126").Value = 50776
$ws.Range("N126").Value = -60656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 84000
$ws.Range("J68").Value = 84000
$ws.Range("L68").Value = 84000
$ws.Range("N68").Value = -85498
$ws.Range("H71").Value = 84000
$ws.Range("J71").Value = 84000
$ws.Range("L71").Value = 252000
$ws.Range("N71").Value = -259488
$ws.Range("H112").Value = 37856.8
$ws.Range("J112").Value = 37856.8
$ws.Range("L112").Value = 37856.8
$ws.Range("N112").Value = -40810.8
$ws.Range("H132").Value = 43494.53
$ws.Range("I132").Value = 1844
$ws.Range("K132").Value = 5532
$ws.Range("M132").Value = -3002
$ws.Range("H134").Value = 2399.925
$ws.Range("I134").Value = 1037.3704
$ws.Range("J134").Value = 5229.846
$ws.Range("K134").Value = 3112.1112
$ws.Range("L134").Value = 15689.538
$ws.Range("M134").Value = -577.1112000000003
$ws.Range("N134").Value = -20759.538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 631.53845
$ws.Range("I6").Value = 110
$ws.Range("J6").Value = 1240
$ws.Range("K6").Value = 330
$ws.Range("L6").Value = 3720
$ws.Range("M6").Value = -217
$ws.Range("N6").Value = -3946
$ws.Range("H131").Value = 910.33
$ws.Range("I131").Value = 722.25
$ws.Range("J131").Value = 918.1667
$ws.Range("K131").Value = 2166.75
$ws.Range("L131").Value = 2754.5001
$ws.Range("M131").Value = 2873.25
$ws.Range("N131").Value = -12834.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 15026.333
$ws.Range("J49").Value = 15026.333
$ws.Range("L49").Value = 15026.333
$ws.Range("N49").Value = -15394.333
$ws.Range("H70").Value = 4517.5
$ws.Range("I70").Value = 4409.304
$ws.Range("J70").Value = 4708.923
$ws.Range("K70").Value = 4409.304
$ws.Range("L70").Value = 4708.923
$ws.Range("M70").Value = -4139.304
$ws.Range("N70").Value = -5248.923
$ws.Range("H73").Value = 4517.5
$ws.Range("I73").Value = 4409.304
$ws.Range("J73").Value = 4708.923
$ws.Range("K73").Value = 4409.304
$ws.Range("L73").Value = 4708.923
$ws.Range("M73").Value = -3473.304
$ws.Range("N73").Value = -6580.923
$ws.Range("H80").Value = 188694.4
$ws.Range("I80").Value = 420974.9
$ws.Range("J80").Value = 2870
$ws.Range("K80").Value = 420974.9
$ws.Range("L80").Value = 2870
$ws.Range("M80").Value = -419976.9
$ws.Range("N80").Value = -4866
$ws.Range("H83").Value = 188694.4
$ws.Range("I83").Value = 420974.9
$ws.Range("J83").Value = 2870
$ws.Range("K83").Value = 2104874.5
$ws.Range("L83").Value = 14350
$ws.Range("M83").Value = -2099882.5
$ws.Range("N83").Value = -24334
$ws.Range("H116").Value = 49515
$ws.Range("J116").Value = 49515
$ws.Range("L116").Value = 49515
$ws.Range("N116").Value = -58693
$ws.Range("H132").Value = 2614.54
$ws.Range("I132").Value = 2226.35
$ws.Range("J132").Value = 4167.3
$ws.Range("K132").Value = 6679.049999999999
$ws.Range("L132").Value = 12501.9
$ws.Range("M132").Value = -4149.049999999999
$ws.Range("N132").Value = -17561.9
$ws.Range("H133").Value = 51520.383
$ws.Range("J133").Value = 51520.383
$ws.Range("L133").Value = 51520.383
$ws.Range("N133").Value = -61640.383

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2114.923
$ws.Range("I40").Value = 1957.9166
$ws.Range("J40").Value = 3999
$ws.Range("K40").Value = 1957.9166
$ws.Range("L40").Value = 3999
$ws.Range("M40").Value = -1821.9166
$ws.Range("N40").Value = -4271
$ws.Range("H123").Value = 28437.6
$ws.Range("J123").Value = 28437.6
$ws.Range("L123").Value = 28437.6
$ws.Range("N123").Value = -38237.6
$ws.Range("H132").Value = 4159.324
$ws.Range("I132").Value = 3841.577
$ws.Range("J132").Value = 4910.364
$ws.Range("K132").Value = 11524.731
$ws.Range("L132").Value = 14731.092
$ws.Range("M132").Value = -8994.731
$ws.Range("N132").Value = -19791.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 63499.5
$ws.Range("J137").Value = 63499.5
$ws.Range("L137").Value = 63499.5
$ws.Range("N137").Value = -73699.5
